# Simulated Wild Card round and logged it
# Appends this game's stats to the 49ers 2021 Team Data workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet - per-play yardage logs (rushing "R" / passing "P") for
# both the offense (OFF) and defense (DEF) get this game's numbers
# appended to the running season log.
# ---------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value2 + " 1 -2 -1 -1 13 2 5 2 16 7 2 4 8 5 6 3 1 9 5 3 12 0 0 5 1 13 8 2 5 1 0 9 1 13 4 4 9 5 -1 6 -2 0 1 5 2 2 -1 7 1 5 11 3 9 26 -1 16 1 13 5 1 5 5 3 0 5 1 2 3 9"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value2 + " 2 1 -1 1 1 4 0 0 -5 -1 2 2 -1 7 2 3 1 18 -2 2 14 1 4 2 1 2 4 -3 4 5 5 8 5 1 5 0 0 1 1 9 3 3 2 4 0 2 5 17"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value2 + " 4 13 -2 6 8 19 -2 8 26 13 31 6 24 32 10 1 21 5 43 14 6 34 9 1 17 11 15 11 11 10 9 19 10 37 5 13 -1 18 5 4"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value2 + " 13 2 46 8 6 7 16 9 7 7 2 15 8 5 24 6 30 4 12 9 10 7 18 20 9 12 8 6 0 5 14 8 6 24 5 6 2 11 10 38 9 10 9"

# ---------------------------------------------------------------------
# OFF sheet - season totals by down/distance, updated with this game's
# offensive production (row 2 = RATT/rushing, row 3 = PATT/passing).
# ---------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("B2").Value = 5
$offWs.Range("C2").Value = 263
$offWs.Range("D2").Value = 15
$offWs.Range("E2").Value = 13
$offWs.Range("F2").Value = 78
$offWs.Range("G2").Value = 69
$offWs.Range("H2").Value = 3
$offWs.Range("I2").Value = 11
$offWs.Range("J2").Value = 39
$offWs.Range("N2").Value = 16

$offWs.Range("C3").Value = 175
$offWs.Range("D3").Value = 3
$offWs.Range("E3").Value = 38
$offWs.Range("F3").Value = 103
$offWs.Range("G3").Value = 29
$offWs.Range("H3").Value = 32
$offWs.Range("I3").Value = 63
$offWs.Range("J3").Value = 46
$offWs.Range("L3").Value = 294
$offWs.Range("M3").Value = 199
$offWs.Range("Q3").Value = 627

# ---------------------------------------------------------------------
# DEF sheet - season totals allowed by down/distance, updated with this
# game's defensive numbers (row 2 = RATT, row 3 = PATT).
# ---------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("C2").Value = 210
$defWs.Range("E2").Value = 16
$defWs.Range("F2").Value = 71
$defWs.Range("G2").Value = 64
$defWs.Range("I2").Value = 7
$defWs.Range("J2").Value = 32
$defWs.Range("N2").Value = 37
$defWs.Range("O2").Value = 27
$defWs.Range("P2").Value = 11

$defWs.Range("B3").Value = 14
$defWs.Range("C3").Value = 204
$defWs.Range("E3").Value = 46
$defWs.Range("F3").Value = 107
$defWs.Range("G3").Value = 35
$defWs.Range("H3").Value = 37
$defWs.Range("I3").Value = 71
$defWs.Range("J3").Value = 68
$defWs.Range("L3").Value = 341
$defWs.Range("M3").Value = 229
$defWs.Range("Q3").Value = 643

# ---------------------------------------------------------------------
# ST sheet - special teams season totals (row 2), touchbacks (row 3),
# plus per-game logs for TB/D/RA/RM (rows 4-6 / col D rows 3-5).
# ---------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 93
$stWs.Range("D2").Value = 59
$stWs.Range("F2").Value = 589
$stWs.Range("G2").Value = 574
$stWs.Range("J2").Value = 290
$stWs.Range("K2").Value = 275
$stWs.Range("L2").Value = 163
$stWs.Range("M2").Value = 128
$stWs.Range("N2").Value = 51
$stWs.Range("O2").Value = 36

$stWs.Range("B3").Value = 41

$stWs.Range("B4").Value = $stWs.Range("B4").Value2 + " 66 69 65 64"
$stWs.Range("B5").Value = $stWs.Range("B5").Value2 + " 14 22 18 32"
$stWs.Range("B6").Value = $stWs.Range("B6").Value2 + " 15 16 13"
$stWs.Range("D3").Value = $stWs.Range("D3").Value2 + " 42 43 47 43 43 46 48 42"
$stWs.Range("D4").Value = $stWs.Range("D4").Value2 + " 9 31 0 0 0 0 5 0"
$stWs.Range("D5").Value = $stWs.Range("D5").Value2 + " 0 0 14 0 11 0 0 0 0"

# ---------------------------------------------------------------------
# TURNS sheet - turnovers forced on the road this game.
# ---------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("B3").Value = 8
$turnsWs.Range("C3").Value = 6
$turnsWs.Range("E3").Value = 15

# ---------------------------------------------------------------------
# PEN sheet - penalty counts.
# ---------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("B2").Value = 23
$penWs.Range("B3").Value = 22
$penWs.Range("D4").Value = 20
